# katalog.xlsx update: add a "Beton" (concrete price) parameter row to the
# Draht_Matten sheet and wire it into the "Montageart" options string so the
# fixed 0 is replaced by the new P_Beton variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Draht_Matten")

# Insert a new row above the old row 7 ("Montage (€/m)" / P_Arbeit), pushing
# the old rows 7-8 down to 8-9.
$ws.Rows.Item(7).Insert()

# Fill in the new row 7: a numeric parameter "Beton" / P_Beton with default 9.
# (Written before the D5 edit below so new shared strings are appended in
# the same order the original authors' edit produced them in.)
$ws.Range("A7").Value = "Zahl"
$ws.Range("B7").Value = "Beton"
$ws.Range("C7").Value = "P_Beton"
$ws.Range("D7").Value = 9

# Update the "Montageart" options (row 5) so "Betonieren" uses the new
# P_Beton variable instead of the hard-coded 0.
$ws.Range("D5").Value = "Betonieren:P_Beton, Konsole:1"

# Match the saved selection/active cell.
$ws.Range("D5").Select()
